$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = '[Alfredo%Guiroy%NULL%1,                    Martín%Gagliardi%NULL%2,                    Martín%Gagliardi%NULL%0,                    Nicolas%Coombes%NULL%1,                    Federico%Landriel%NULL%2,                    Federico%Landriel%NULL%0,                    Carlos%Zanardi%NULL%1,                    Gastón Camino%Willhuber%NULL%1,                    Juan Pablo%Guyot%NULL%1,                    Marcelo%Valacco%NULL%1]'
$ws.Cells.Item(3, 5).Value = '[Yufei%Li%NULL%10,                    Nathaniel%Scherer%NULL%10,                    Lambert%Felix%NULL%20,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%20,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%20,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(5, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(6, 5).Value = '[Kossi Blewussi%Kounou%benkounou@hotmail.fr%1,                    Koffi Mawuse%Guédénon%NULL%2,                    Koffi Mawuse%Guédénon%NULL%0,                    Ayoko Akouavi%Dogbe Foli%NULL%1,                    Eric%Gnassounou‐Akpa%NULL%1]'
$ws.Cells.Item(7, 5).Value = '[Jianbo%Lai%NULL%0,                    Simeng%Ma%NULL%0,                    Ying%Wang%NULL%0,                    Zhongxiang%Cai%NULL%0,                    Jianbo%Hu%NULL%0,                    Ning%Wei%NULL%0,                    Jiang%Wu%NULL%0,                    Hui%Du%NULL%0,                    Tingting%Chen%NULL%0,                    Ruiting%Li%NULL%0,                    Huawei%Tan%NULL%0,                    Lijun%Kang%NULL%0,                    Lihua%Yao%NULL%0,                    Manli%Huang%NULL%0,                    Huafen%Wang%NULL%0,                    Gaohua%Wang%NULL%0,                    Zhongchun%Liu%NULL%0,                    Shaohua%Hu%NULL%0]'
$ws.Cells.Item(8, 5).Value = '[Simon Ching%Lam%NULL%1,                    Teresa%Arora%NULL%1,                    Ian%Grey%NULL%1,                    Lorna Kwai Ping%Suen%NULL%1,                    Emma Yun-zhi%Huang%NULL%1,                    Daofan%Li%NULL%1,                    Kin Bong Hubert%Lam%NULL%1]'
$ws.Cells.Item(9, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(10, 5).Value = '[Ruilin%Li%NULL%1,                    Youlin%Chen%NULL%1,                    Jianlin%Lv%NULL%1,                    Linlin%Liu%NULL%1,                    Shiqin%Zong%NULL%1,                    Hanxia%Li%NULL%1,                    Hong%Li%NULL%1,                    Massimo%Tusconi.%NULL%2,                    Massimo%Tusconi.%NULL%0]'
$ws.Cells.Item(11, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(12, 5).Value = '[Wen%Lu%luwen67@sina.com%0,                    Hang%Wang%NULL%1,                    Yuxing%Lin%NULL%1,                    Li%Li%lilifuzhou@126.com%3]'
$ws.Cells.Item(13, 5).Value = '[Abdallah Y.%Naser%abdallah.naser@iu.edu.jo%0,                    Eman Zmaily%Dahmash%NULL%2,                    Eman Zmaily%Dahmash%NULL%0,                    Rabaa%Al‐Rousan%NULL%2,                    Rabaa%Al‐Rousan%NULL%0,                    Hassan%Alwafi%NULL%1,                    Hamzeh Mohammad%Alrawashdeh%NULL%1,                    Imene%Ghoul%NULL%1,                    Anwer%Abidine%NULL%1,                    Mohammed A.%Bokhary%NULL%1,                    Hadeel T.%AL‐Hadithi%NULL%1,                    Dalia%Ali%NULL%1,                    Rasha%Abuthawabeh%NULL%1,                    Ghada Mohammad%Abdelwahab%NULL%1,                    Yosra J.%Alhartani%NULL%1,                    Haneen%Al Muhaisen%NULL%1,                    Ayah%Dagash%NULL%1,                    Hamad S.%Alyami%NULL%1]'
$ws.Cells.Item(14, 5).Value = '[Moluk%Pouralizadeh%NULL%1,                    Zahra%Bostani%NULL%1,                    Saman%Maroufizadeh%NULL%1,                    Atefeh%Ghanbari%NULL%1,                    Maryam%Khoshbakht%NULL%1,                    Seyed Amirhossein%Alavi%NULL%1,                    Sadra%Ashrafi%NULL%1]'
$ws.Cells.Item(16, 5).Value = '[Jianyu%Que%NULL%1,                    Le%Shi%NULL%2,                    Le%Shi%NULL%0,                    Jiahui%Deng%NULL%1,                    Jiajia%Liu%NULL%1,                    Li%Zhang%NULL%0,                    Suying%Wu%NULL%1,                    Yimiao%Gong%NULL%1,                    Weizhen%Huang%NULL%1,                    Kai%Yuan%NULL%1,                    Wei%Yan%NULL%1,                    Yankun%Sun%NULL%1,                    Maosheng%Ran%NULL%1,                    Yanping%Bao%NULL%1,                    Lin%Lu%NULL%1]'
$ws.Cells.Item(17, 5).Value = '[Rodolfo%Rossi%NULL%0,                    Valentina%Socci%NULL%0,                    Francesca%Pacitti%NULL%0,                    Giorgio%Di Lorenzo%NULL%0,                    Antinisca%Di Marco%NULL%0,                    Alberto%Siracusano%NULL%0,                    Alessandro%Rossi%NULL%0]'
$ws.Cells.Item(18, 5).Value = '[ M.%Salman%null%1,                     M. H.% Raza%null%1,                     Z. U.% Mustafa%null%1,                     T. M.% Khan%null%1,                     N.% Asif%null%1,                     H.% Tahir%null%1,                     N.% Shehzadi%null%1,                     K. % Hussain%null%1]'
$ws.Cells.Item(19, 5).Value = '[Ari%Shechter%NULL%1,                    Franchesca%Diaz%NULL%1,                    Nathalie%Moise%NULL%1,                    D. Edmund%Anstey%NULL%1,                    Siqin%Ye%NULL%1,                    Sachin%Agarwal%NULL%1,                    Jeffrey L.%Birk%NULL%1,                    Daniel%Brodie%NULL%0,                    Diane E.%Cannone%NULL%1,                    Bernard%Chang%NULL%1,                    Jan%Claassen%NULL%2,                    Talea%Cornelius%NULL%1,                    Lilly%Derby%NULL%1,                    Melissa%Dong%NULL%1,                    Raymond C.%Givens%NULL%1,                    Beth%Hochman%NULL%1,                    Shunichi%Homma%NULL%1,                    Ian M.%Kronish%NULL%1,                    Sung A.J.%Lee%NULL%1,                    Wilhelmina%Manzano%NULL%1,                    Laurel E.S.%Mayer%NULL%1,                    Cara L.%McMurry%NULL%1,                    Vivek%Moitra%NULL%1,                    Patrick%Pham%NULL%1,                    LeRoy%Rabbani%NULL%1,                    Reynaldo R.%Rivera%NULL%1,                    Allan%Schwartz%NULL%2,                    Joseph E.%Schwartz%NULL%1,                    Peter A.%Shapiro%NULL%1,                    Kaitlin%Shaw%NULL%1,                    Alexandra M.%Sullivan%NULL%1,                    Courtney%Vose%NULL%1,                    Lauren%Wasson%NULL%1,                    Donald%Edmondson%NULL%1,                    Marwah%Abdalla%NULL%1]'
$ws.Cells.Item(20, 5).Value = '[Xingyue%Song%NULL%0,                    Wenning%Fu%NULL%0,                    Xiaoran%Liu%NULL%0,                    Zhiqian%Luo%NULL%0,                    Rixing%Wang%NULL%0,                    Ning%Zhou%NULL%0,                    Shijiao%Yan%NULL%0,                    Chuanzhu%Lv%NULL%0]'
$ws.Cells.Item(21, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(22, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(23, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(24, 5).Value = '[Ya-Xi%Wang%NULL%1,                    Hong-Tao%Guo%NULL%1,                    Xue-Wei%Du%NULL%1,                    Wen%Song%NULL%1,                    Chang%Lu%NULL%1,                    Wen-Nv%Hao%NULL%1,                    Maria%Kapritsou.%NULL%3,                    Maria%Kapritsou.%NULL%0,                    Maria%Kapritsou.%NULL%0]'
$ws.Cells.Item(25, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(26, 5).Value = '[Xiao%Xiao%NULL%0,                    Xiaobin%Zhu%NULL%1,                    Shuai%Fu%NULL%1,                    Yugang%Hu%NULL%1,                    Xiaoning%Li%NULL%1,                    Jinsong%Xiao%NULL%1]'
$ws.Cells.Item(27, 5).Value = '[Stephen X.%Zhang%NULL%0,                    Jing%Liu%NULL%0,                    Asghar%Afshar Jahanshahi%NULL%1,                    Khaled%Nawaser%NULL%1,                    Ali%Yousefi%NULL%1,                    Jizhen%Li%NULL%1,                    Shuhua%Sun%NULL%1]'
$ws.Cells.Item(28, 5).Value = '[Wen-rui%Zhang%NULL%0,                    Kun%Wang%NULL%2,                    Lu%Yin%NULL%0,                    Wen-feng%Zhao%NULL%1,                    Qing%Xue%NULL%1,                    Mao%Peng%NULL%1,                    Bao-quan%Min%NULL%1,                    Qing%Tian%NULL%1,                    Hai-xia%Leng%NULL%1,                    Jia-lin%Du%NULL%1,                    Hong%Chang%NULL%1,                    Yuan%Yang%NULL%0,                    Wei%Li%NULL%2,                    Fang-fang%Shangguan%NULL%1,                    Tian-yi%Yan%NULL%1,                    Hui-qing%Dong%NULL%1,                    Ying%Han%NULL%1,                    Yu-ping%Wang%NULL%1,                    Fiammetta%Cosci%NULL%1,                    Hong-xing%Wang%NULL%1]'
$ws.Cells.Item(29, 5).Value = '[Yeen%Huang%NULL%0,                    Ning%Zhao%zhaoning2018@email.szu.edu.cn%1]'
$ws.Cells.Item(30, 5).Value = '[Juhong%Zhu%NULL%0,                    Lin%Sun%NULL%1,                    Lan%Zhang%NULL%1,                    Huan%Wang%NULL%1,                    Ajiao%Fan%NULL%1,                    Bin%Yang%NULL%1,                    Wei%Li%NULL%0,                    Shifu%Xiao%NULL%1]'
$ws.Cells.Item(31, 5).Value = '[Jiang%Du%NULL%0,                    Lu%Dong%NULL%1,                    Tao%Wang%NULL%0,                    Chenxin%Yuan%NULL%1,                    Rao%Fu%NULL%1,                    Lei%Zhang%NULL%2,                    Bo%Liu%NULL%1,                    Mingmin%Zhang%NULL%1,                    Yuanyuan%Yin%NULL%1,                    Jiawen%Qin%NULL%1,                    Jennifer%Bouey%NULL%1,                    Min%Zhao%NULL%1,                    Xin%Li%NULL%2]'
$ws.Cells.Item(32, 5).Value = '[Rümeysa Yeni%Elbay%rumeysa.yenielbay@medeniyet.edu.tr%0,                    Ayşe%Kurtulmuş%NULL%1,                    Selim%Arpacıoğlu%NULL%1,                    Emrah%Karadere%NULL%1]'
$ws.Cells.Item(33, 5).Value = '[Anne%Moorhead%NULL%1,                    Gunther%Eysenbach%NULL%1,                    Albert%Wu%NULL%2,                    Albert%Wu%NULL%0,                    Bojana%Knezevic%NULL%1,                    Bradley A%Evanoff%bevanoff@wustl.edu%2,                    Bradley A%Evanoff%bevanoff@wustl.edu%0,                    Jaime R%Strickland%NULL%2,                    Jaime R%Strickland%NULL%0,                    Ann Marie%Dale%NULL%2,                    Ann Marie%Dale%NULL%0,                    Lisa%Hayibor%NULL%2,                    Lisa%Hayibor%NULL%0,                    Emily%Page%NULL%2,                    Emily%Page%NULL%0,                    Jennifer G%Duncan%NULL%2,                    Jennifer G%Duncan%NULL%0,                    Thomas%Kannampallil%NULL%2,                    Thomas%Kannampallil%NULL%0,                    Diana L%Gray%NULL%2,                    Diana L%Gray%NULL%0]'
$ws.Cells.Item(34, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(35, 5).Value = '[Seshadri Sekhar%Chatterjee%NULL%0,                    Ranjan%Bhattacharyya%NULL%1,                    Sumita%Bhattacharyya%NULL%1,                    Sukanya%Gupta%NULL%1,                    Soumitra%Das%NULL%1,                    Bejoy Bikram%Banerjee%NULL%1]'
$ws.Cells.Item(36, 5).Value = '[Jie%Chen%NULL%0,                    Xinghuang%Liu%NULL%2,                    Xinghuang%Liu%NULL%0,                    Dongke%Wang%NULL%2,                    Dongke%Wang%NULL%0,                    Yan%Jin%NULL%2,                    Yan%Jin%NULL%0,                    Miao%He%NULL%1,                    Yanling%Ma%NULL%1,                    Xiaolong%Zhao%NULL%1,                    Shuangning%Song%NULL%1,                    Lei%Zhang%NULL%0,                    Xuelian%Xiang%NULL%1,                    Ling%Yang%NULL%0,                    Jun%Song%song111jun@126.com%1,                    Tao%Bai%drbaitao@126.com%4,                    Xiaohua%Hou%houxh@hust.edu.cn%2,                    Xiaohua%Hou%houxh@hust.edu.cn%0]'
$ws.Cells.Item(37, 5).Value = '[Yun%Chen%NULL%0,                    Hao%Zhou%haoye320@163.com%1,                    Yan%Zhou%NULL%1,                    Fang%Zhou%NULL%1]'
$ws.Cells.Item(38, 5).Value = '[Tawfiq%Choudhury%xref no email%1,            Maciej%Debski%xref no email%1,            Andrew%Wiper%xref no email%1,            Amr%Abdelrahman%xref no email%1,            Susan%Wild%xref no email%1,            Shajil%Chalil%xref no email%1,            Ranjit%More%xref no email%1,            Grahame%Goode%xref no email%1,            Billal%Patel%xref no email%1,            Hesham K.%Abdelaziz%xref no email%1]'
$ws.Cells.Item(39, 5).Value = '[Anucha%Apisarnthanarak%NULL%1,                    Piyaporn%Apisarnthanarak%NULL%2,                    Piyaporn%Apisarnthanarak%NULL%0,                    Chanida%Siripraparat%NULL%2,                    Chanida%Siripraparat%NULL%0,                    Pavarat%Saengaram%NULL%1,                    Narakorn%Leeprechanon%NULL%1,                    David J.%Weber%NULL%1]'
$ws.Cells.Item(40, 5).Value = '[Andrea%Amerio%NULL%1,                    Davide%Bianchi%NULL%1,                    Francesca%Santi%NULL%1,                    Luigi%Costantini%NULL%1,                    Anna%Odone%NULL%1,                    Carlo%Signorelli%NULL%0,                    Alessandra%Costanza%NULL%1,                    Gianluca%Serafini%NULL%1,                    Mario%Amore%NULL%1,                    Andrea%Aguglia%NULL%1]'
$ws.Cells.Item(41, 5).Value = '[Yufei%Li%NULL%0,                    Nathaniel%Scherer%NULL%0,                    Lambert%Felix%NULL%0,                    Lambert%Felix%NULL%0,                    Hannah%Kuper%NULL%0,                    Hannah%Kuper%NULL%0,                    Jakob%Pietschnig%NULL%0,                    Jakob%Pietschnig%NULL%0]'
$ws.Cells.Item(2, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(3, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(5, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(6, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(7, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(8, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(9, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(10, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(11, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(12, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(13, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(14, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(16, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(17, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(18, 9).Value = '_MedBiorxiv_CROSSREF'
$ws.Cells.Item(19, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(20, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(21, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(22, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(23, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(24, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(25, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(26, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(27, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(28, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(29, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(30, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(31, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(32, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(33, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(34, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(35, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(36, 9).Value = '_PMC_Springer_CROSSREF'
$ws.Cells.Item(37, 9).Value = '_PMC_elsevier_CROSSREF'
$ws.Cells.Item(38, 9).Value = '_CROSSREF'
$ws.Cells.Item(39, 9).Value = '_PMC_CROSSREF'
$ws.Cells.Item(41, 9).Value = '_PMC_CROSSREF'
